$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(233, 3).Value = 1
$ws.Cells.Item(233, 4).Value = 9958.450821037621
$ws.Cells.Item(233, 5).Value = 11664
$ws.Cells.Item(233, 6).Value = 984135.6916657642
$ws.Cells.Item(233, 8).Value = 654.456344506836
$ws.Cells.Item(233, 9).Value = 974831.6971892334

$ws.Cells.Item(234, 4).Value = 981965.8571301685
$ws.Cells.Item(234, 6).Value = 981965.8571301685
$ws.Cells.Item(234, 8).Value = 653.5465717285157
$ws.Cells.Item(234, 10).Value = -2824.290880102548

$ws.Cells.Item(235, 4).Value = 981965.8571301685
$ws.Cells.Item(235, 6).Value = 981965.8571301685

$ws.Cells.Item(236, 4).Value = 981965.8571301685
$ws.Cells.Item(236, 6).Value = 981965.8571301685

$ws.Cells.Item(237, 4).Value = 981965.8571301685
$ws.Cells.Item(237, 6).Value = 981965.8571301685

$ws.Cells.Item(238, 4).Value = 981965.8571301685
$ws.Cells.Item(238, 6).Value = 981965.8571301685

$ws.Cells.Item(239, 4).Value = 981965.8571301685
$ws.Cells.Item(239, 6).Value = 981965.8571301685

$ws.Cells.Item(240, 4).Value = 981965.8571301685
$ws.Cells.Item(240, 6).Value = 981965.8571301685

$ws.Cells.Item(241, 4).Value = 981965.8571301685
$ws.Cells.Item(241, 6).Value = 981965.8571301685

$ws.Cells.Item(242, 4).Value = 9960.877026016555
$ws.Cells.Item(242, 5).Value = 12263
$ws.Cells.Item(242, 6).Value = 981313.0957989048
$ws.Cells.Item(242, 8).Value = 652.7613312637329
$ws.Cells.Item(242, 9).Value = 972004.9801041519

$ws.Cells.Item(243, 4).Value = 976983.651157954
$ws.Cells.Item(243, 6).Value = 976983.651157954
$ws.Cells.Item(243, 8).Value = 650.5539968772889
$ws.Cells.Item(243, 10).Value = -4982.205972214462

$ws.Cells.Item(244, 4).Value = 976983.651157954
$ws.Cells.Item(244, 6).Value = 976983.651157954

$ws.Cells.Item(245, 4).Value = 976983.651157954
$ws.Cells.Item(245, 6).Value = 976983.651157954

$ws.Cells.Item(246, 4).Value = 976983.651157954
$ws.Cells.Item(246, 6).Value = 976983.651157954

$ws.Cells.Item(247, 4).Value = 976983.651157954
$ws.Cells.Item(247, 6).Value = 976983.651157954

$ws.Cells.Item(248, 4).Value = 9891.165019591124
$ws.Cells.Item(248, 5).Value = 12235
$ws.Cells.Item(248, 6).Value = 976333.8355556324
$ws.Cells.Item(248, 8).Value = 649.8156023216249
$ws.Cells.Item(248, 9).Value = 967092.4861383629

$ws.Cells.Item(249, 4).Value = 983387.4581838596
$ws.Cells.Item(249, 6).Value = 983387.4581838596
$ws.Cells.Item(249, 8).Value = 654.4404401626588
$ws.Cells.Item(249, 10).Value = 6403.807025905582

$ws.Cells.Item(250, 4).Value = 983387.4581838596
$ws.Cells.Item(250, 6).Value = 983387.4581838596

$ws.Cells.Item(251, 4).Value = 983387.4581838596
$ws.Cells.Item(251, 6).Value = 983387.4581838596

$ws.Cells.Item(252, 4).Value = 10012.57356101871
$ws.Cells.Item(252, 5).Value = 12070
$ws.Cells.Item(252, 6).Value = 982733.8754027545
$ws.Cells.Item(252, 8).Value = 653.5827811050416
$ws.Cells.Item(252, 9).Value = 973374.8846228409

$ws.Cells.Item(253, 4).Value = 988232.2789547107
$ws.Cells.Item(253, 6).Value = 988232.2789547107
$ws.Cells.Item(253, 8).Value = 657.276188949585
$ws.Cells.Item(253, 10).Value = 4844.820770851104

$ws.Cells.Item(254, 4).Value = 988232.2789547107
$ws.Cells.Item(254, 6).Value = 988232.2789547107

$ws.Cells.Item(255, 4).Value = 10064.79339142029
$ws.Cells.Item(255, 5).Value = 12074
$ws.Cells.Item(255, 6).Value = 987575.8223373431
$ws.Cells.Item(255, 8).Value = 656.4566173675538
$ws.Cells.Item(255, 9).Value = 978167.4855632904
